$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bench Ratio row for all four days, entered first -------------------
# (D6 keeps its existing shared-string slot; that string's text becomes
# "Bench Ratio (Y:N): 12:5" simply by overwriting it in place.)
$ws.Range("D6").Value  = "Bench Ratio (Y:N): 12:5"
$ws.Range("H6").Value  = "Bench Ratio (Y:N): 13:7"
$ws.Range("L6").Value  = "Bench Ratio (Y:N): 12:16"
$ws.Range("P6").Value  = "Bench Ratio (Y:N): 28:17"

# --- Monday block (column D): totals, then averages -----------------------
$ws.Range("D9").Value  = "Total Wait: 234 min"
$ws.Range("D10").Value = "Total Lift: 637 min"
$ws.Range("D7").Value  = "Average Wait: 7.0909"
$ws.Range("D8").Value  = "Average Lift: 28.9545"

# --- Tuesday block (column H): totals, then averages -----------------------
$ws.Range("H9").Value  = "Total Wait: 112"
$ws.Range("H7").Value  = "Average Wait: 4.8700"
$ws.Range("H10").Value = "Total Lift: 352 "
$ws.Range("H8").Value  = "Average Lift: 23.4667"

# --- Wednesday block (column L): totals, then averages ----------------------
$ws.Range("L9").Value  = "Total Wait: 156"
$ws.Range("L7").Value  = "Average Wait: 3.8049"
$ws.Range("L10").Value = "Total Lift: 715"
$ws.Range("L8").Value  = "Average Lift: 22.3438"

# --- Thursday block (column P): totals, then averages -----------------------
$ws.Range("P9").Value  = "Total Wait: 150"
$ws.Range("P10").Value = "Total Lift: 989"
$ws.Range("P7").Value  = "Average Wait: 3.125"
$ws.Range("P8").Value  = "Average Lift: 24.1220"

# --- Remove now-superfluous blank styled cells ---------------------------
$ws.Range("G28").Clear()
$ws.Range("G29").Clear()
$ws.Range("G30").Clear()
$ws.Range("C36").Clear()
$ws.Range("C37").Clear()

# --- Add new blank styled cells ------------------------------------------
$ws.Range("A38").HorizontalAlignment = -4108

$ws.Range("I44:J44").Font.Bold = $true
$ws.Range("M51:N51").Font.Bold = $true

# --- Selection / view -----------------------------------------------------
$ws.Range("I44:J44").Select()
